# Updated cryptos list data (Price and Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "30.589.78"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.875.80"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2901"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07734"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.875.55"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "30.596.37"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007470"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "2.123.36"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.213"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.186"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09881"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.346"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.508"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.243"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04767"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6935"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.223"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.968"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4160"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8343"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.392"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "920.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05671"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.73%  "
